$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value2 = 1139.3
$ws.Range("J19").Value2 = 1154.6666
$ws.Range("L19").Value2 = 1154.6666
$ws.Range("N19").Value2 = -1504.6666
$ws.Range("H33").Value2 = 773.2105
$ws.Range("I33").Value2 = 847.5333000000001
$ws.Range("K33").Value2 = 847.5333000000001
$ws.Range("M33").Value2 = -618.5333000000001
$ws.Range("H92").Value2 = 142857800
$ws.Range("I92").Value2 = 712.6
$ws.Range("J92").Value2 = 500000480
$ws.Range("K92").Value2 = 712.6
$ws.Range("L92").Value2 = 500000480
$ws.Range("M92").Value2 = 535.4
$ws.Range("N92").Value2 = -500002976
$ws.Range("H132").Value2 = 2964
$ws.Range("I132").Value2 = 3014.5
$ws.Range("K132").Value2 = 9043.5
$ws.Range("M132").Value2 = -6513.5
$ws.Range("H137").Value2 = 2486.2778
$ws.Range("I137").Value2 = 2689.0833
$ws.Range("K137").Value2 = 8067.249899999999
$ws.Range("M137").Value2 = -5517.249899999999
$ws.Range("H138").Value2 = 3613.18
$ws.Range("J138").Value2 = 5768.9644
$ws.Range("L138").Value2 = 17306.8932
$ws.Range("N138").Value2 = -27586.8932

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value2 = 0
$ws.Range("J23").Value2 = 0
$ws.Range("L23").Value2 = 0
$ws.Range("N23").ClearContents()
$ws.Range("H32").Value2 = 1895961.2
$ws.Range("I32").Value2 = 1955046.8
$ws.Range("K32").Value2 = 1955046.8
$ws.Range("M32").Value2 = -1954759.8
$ws.Range("H45").Value2 = 3691.0908
$ws.Range("I45").Value2 = 1769.75
$ws.Range("K45").Value2 = 1769.75
$ws.Range("M45").Value2 = -1392.75
$ws.Range("H61").Value2 = 2186.8484
$ws.Range("I61").Value2 = 2208.6128
$ws.Range("K61").Value2 = 2208.6128
$ws.Range("M61").Value2 = -1996.6128
$ws.Range("H74").Value2 = 48771.17
$ws.Range("I74").Value2 = 66611.28
$ws.Range("J74").Value2 = 4170.9
$ws.Range("K74").Value2 = 66611.28
$ws.Range("L74").Value2 = 4170.9
$ws.Range("M74").Value2 = -65737.28
$ws.Range("N74").Value2 = -5918.9
$ws.Range("H77").Value2 = 48771.17
$ws.Range("I77").Value2 = 66611.28
$ws.Range("J77").Value2 = 4170.9
$ws.Range("K77").Value2 = 333056.4
$ws.Range("L77").Value2 = 20854.5
$ws.Range("M77").Value2 = -328688.4
$ws.Range("N77").Value2 = -29590.5
$ws.Range("H122").Value2 = 200000
$ws.Range("I122").Value2 = 200000
$ws.Range("J122").Value2 = 0
$ws.Range("K122").Value2 = 600000
$ws.Range("L122").Value2 = 0
$ws.Range("M122").Value2 = -597550
$ws.Range("N122").ClearContents()
$ws.Range("H136").Value2 = 2186.8484
$ws.Range("I136").Value2 = 2208.6128
$ws.Range("K136").Value2 = 6625.8384
$ws.Range("M136").Value2 = -4075.8384

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value2 = 3639888.8
$ws.Range("J99").Value2 = 6997446
$ws.Range("L99").Value2 = 6997446
$ws.Range("N99").Value2 = -7000442

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value2 = 280.5263
$ws.Range("I7").Value2 = 211.33333
$ws.Range("K7").Value2 = 211.33333
$ws.Range("M7").Value2 = -98.33332999999999
$ws.Range("H16").Value2 = 6183.3687
$ws.Range("I16").Value2 = 3786.8572
$ws.Range("K16").Value2 = 3786.8572
$ws.Range("M16").Value2 = -3499.8572
$ws.Range("H31").Value2 = 9168.578
$ws.Range("I31").Value2 = 3859
$ws.Range("J31").Value2 = 11325.594
$ws.Range("K31").Value2 = 3859
$ws.Range("L31").Value2 = 11325.594
$ws.Range("M31").Value2 = -3564
$ws.Range("N31").Value2 = -11915.594
$ws.Range("H34").Value2 = 9168.578
$ws.Range("I34").Value2 = 3859
$ws.Range("J34").Value2 = 11325.594
$ws.Range("K34").Value2 = 3859
$ws.Range("L34").Value2 = 11325.594
$ws.Range("M34").Value2 = -3657
$ws.Range("N34").Value2 = -11729.594
$ws.Range("H88").Value2 = 40343
$ws.Range("J88").Value2 = 40343
$ws.Range("L88").Value2 = 40343
$ws.Range("N88").Value2 = -41155
$ws.Range("H91").Value2 = 40343
$ws.Range("J91").Value2 = 40343
$ws.Range("L91").Value2 = 40343
$ws.Range("N91").Value2 = -43151
$ws.Range("H96").Value2 = 47634.57
$ws.Range("J96").Value2 = 47634.57
$ws.Range("L96").Value2 = 47634.57
$ws.Range("N96").Value2 = -53126.57
$ws.Range("H105").Value2 = 4188.9287
$ws.Range("I105").Value2 = 1738.3334
$ws.Range("K105").Value2 = 1738.3334
$ws.Range("M105").Value2 = 8.666600000000017
$ws.Range("H111").Value2 = 0
$ws.Range("J111").Value2 = 0
$ws.Range("L111").Value2 = 0
$ws.Range("N111").ClearContents()
$ws.Range("H113").Value2 = 6183.3687
$ws.Range("I113").Value2 = 3786.8572
$ws.Range("K113").Value2 = 3786.8572
$ws.Range("M113").Value2 = -1616.8572
$ws.Range("H132").Value2 = 6943
$ws.Range("I132").Value2 = 2558.2222
$ws.Range("K132").Value2 = 7674.6666
$ws.Range("M132").Value2 = -5144.6666
$ws.Range("H134").Value2 = 8954.796
$ws.Range("I134").Value2 = 9765.315000000001
$ws.Range("J134").Value2 = 8338.799999999999
$ws.Range("K134").Value2 = 29295.945
$ws.Range("L134").Value2 = 25016.4
$ws.Range("M134").Value2 = -26760.945
$ws.Range("N134").Value2 = -30086.4
$ws.Range("H139").Value2 = 62799.4
$ws.Range("J139").Value2 = 67249.25
$ws.Range("L139").Value2 = 67249.25
$ws.Range("N139").Value2 = -77529.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I2").Value2 = 373.75
$ws.Range("J2").Value2 = 572911.3
$ws.Range("K2").Value2 = 2242.5
$ws.Range("L2").Value2 = 3437467.8
$ws.Range("M2").Value2 = -2129.5
$ws.Range("N2").Value2 = -3437693.8
$ws.Range("H12").Value2 = 2500736.2
$ws.Range("J12").Value2 = 3572132.8
$ws.Range("L12").Value2 = 10716398.4
$ws.Range("N12").Value2 = -10716744.4
$ws.Range("H40").Value2 = 142.75
$ws.Range("I40").Value2 = 85.5
$ws.Range("J40").Value2 = 200
$ws.Range("K40").Value2 = 342
$ws.Range("L40").Value2 = 800
$ws.Range("M40").Value2 = -273
$ws.Range("N40").Value2 = -938
$ws.Range("H82").Value2 = 19000
$ws.Range("I82").Value2 = 19000
$ws.Range("K82").Value2 = 57000
$ws.Range("M82").Value2 = -56594
$ws.Range("H85").Value2 = 19000
$ws.Range("I85").Value2 = 19000
$ws.Range("K85").Value2 = 57000
$ws.Range("M85").Value2 = -55596

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value2 = 1250343.2
$ws.Range("I2").Value2 = 94.333336
$ws.Range("J2").Value2 = 2857806.2
$ws.Range("K2").Value2 = 94.333336
$ws.Range("L2").Value2 = 2857806.2
$ws.Range("M2").Value2 = 18.666664
$ws.Range("N2").Value2 = -2858032.2
$ws.Range("H122").Value2 = 999999
$ws.Range("I122").Value2 = 999999
$ws.Range("J122").Value2 = 0
$ws.Range("K122").Value2 = 2999997
$ws.Range("L122").Value2 = 0
$ws.Range("M122").Value2 = -2997547
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value2 = 9096.5
$ws.Range("I132").Value2 = 4674.6665
$ws.Range("J132").Value2 = 13518.333
$ws.Range("K132").Value2 = 14023.9995
$ws.Range("L132").Value2 = 40554.999
$ws.Range("M132").Value2 = -11493.9995
$ws.Range("N132").Value2 = -45614.999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 6670.4346
$ws.Range("I7").Value2 = 5117.9165
$ws.Range("J7").Value2 = 8364.091
$ws.Range("K7").Value2 = 5117.9165
$ws.Range("L7").Value2 = 8364.091
$ws.Range("M7").Value2 = -5005.9165
$ws.Range("N7").Value2 = -8588.091
$ws.Range("H40").Value2 = 8311.75
$ws.Range("I40").Value2 = 7559.2
$ws.Range("J40").Value2 = 8849.286
$ws.Range("K40").Value2 = 7559.2
$ws.Range("L40").Value2 = 8849.286
$ws.Range("M40").Value2 = -7423.2
$ws.Range("N40").Value2 = -9121.286
$ws.Range("H61").Value2 = 7828.5713
$ws.Range("I61").Value2 = 6500
$ws.Range("K61").Value2 = 6500
$ws.Range("M61").Value2 = -6298
$ws.Range("H87").Value2 = 0
$ws.Range("J87").Value2 = 0
$ws.Range("L87").Value2 = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value2 = 0
$ws.Range("J90").Value2 = 0
$ws.Range("L90").Value2 = 0
$ws.Range("N90").ClearContents()
$ws.Range("H93").Value2 = 7200.5
$ws.Range("I93").Value2 = 6640.6
$ws.Range("K93").Value2 = 6640.6
$ws.Range("M93").Value2 = -5392.6
$ws.Range("H113").Value2 = 7828.5713
$ws.Range("I113").Value2 = 6500
$ws.Range("K113").Value2 = 6500
$ws.Range("M113").Value2 = -4330
$ws.Range("H122").Value2 = 3601.0227
$ws.Range("I122").Value2 = 2935.4375
$ws.Range("J122").Value2 = 5375.9165
$ws.Range("K122").Value2 = 8806.3125
$ws.Range("L122").Value2 = 16127.7495
$ws.Range("M122").Value2 = -6356.3125
$ws.Range("N122").Value2 = -21027.7495
$ws.Range("H126").Value2 = 6670.4346
$ws.Range("I126").Value2 = 5117.9165
$ws.Range("J126").Value2 = 8364.091
$ws.Range("K126").Value2 = 15353.7495
$ws.Range("L126").Value2 = 25092.273
$ws.Range("M126").Value2 = -12883.7495
$ws.Range("N126").Value2 = -30032.273
$ws.Range("H136").Value2 = 8401.206
$ws.Range("J136").Value2 = 11029
$ws.Range("L136").Value2 = 33087
$ws.Range("N136").Value2 = -38187

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value2 = 18519240
$ws.Range("I107").Value2 = 812.8570999999999
$ws.Range("J107").Value2 = 30303692
$ws.Range("K107").Value2 = 2438.5713
$ws.Range("L107").Value2 = 90911076
$ws.Range("M107").Value2 = -518.5712999999996
$ws.Range("N107").Value2 = -90914916
$ws.Range("H113").Value2 = 2944.3333
$ws.Range("J113").Value2 = 2944.3333
$ws.Range("L113").Value2 = 8832.999899999999
$ws.Range("N113").Value2 = -13172.9999
$ws.Range("H126").Value2 = 5710.1
$ws.Range("I126").Value2 = 4580.6
$ws.Range("K126").Value2 = 13741.8
$ws.Range("M126").Value2 = -11271.8
